$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1909.0454
$ws.Range("I17").Value = 1300
$ws.Range("J17").Value = 2005.2106
$ws.Range("K17").Value = 3900
$ws.Range("L17").Value = 6015.6318
$ws.Range("M17").Value = -3732
$ws.Range("N17").Value = -6351.6318
$ws.Range("H40").Value = 1500
$ws.Range("I40").Value = 1500
$ws.Range("K40").Value = 1500
$ws.Range("M40").Value = -1325
$ws.Range("H92").Value = 994.1905
$ws.Range("I92").Value = 257.52942
$ws.Range("K92").Value = 257.52942
$ws.Range("M92").Value = 990.4705799999999
$ws.Range("H98").Value = 1426.5264
$ws.Range("I98").Value = 1502.9375
$ws.Range("K98").Value = 1502.9375
$ws.Range("M98").Value = -4.9375
$ws.Range("H100").Value = 2010.5555
$ws.Range("I100").Value = 1973.75
$ws.Range("J100").Value = 2084.1667
$ws.Range("K100").Value = 1973.75
$ws.Range("L100").Value = 2084.1667
$ws.Range("M100").Value = -1432.75
$ws.Range("N100").Value = -3166.1667
$ws.Range("H122").Value = 1426.5264
$ws.Range("I122").Value = 1502.9375
$ws.Range("K122").Value = 4508.8125
$ws.Range("M122").Value = -2058.8125
$ws.Range("H132").Value = 2410.7
$ws.Range("I132").Value = 2274.4211
$ws.Range("K132").Value = 6823.263300000001
$ws.Range("M132").Value = -4293.263300000001
$ws.Range("H141").Value = 8832.52
$ws.Range("I141").Value = 7290.65
$ws.Range("K141").Value = 21871.95
$ws.Range("M141").Value = -16691.95

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9354.656999999999
$ws.Range("I32").Value = 4355.829
$ws.Range("K32").Value = 4355.829
$ws.Range("M32").Value = -4068.829
$ws.Range("H110").Value = 1549.75
$ws.Range("I110").Value = 1099.5
$ws.Range("K110").Value = 1099.5
$ws.Range("M110").Value = 945.5
$ws.Range("H122").Value = 2151.25
$ws.Range("I122").Value = 1785
$ws.Range("K122").Value = 5355
$ws.Range("M122").Value = -2905
$ws.Range("H132").Value = 643397.1
$ws.Range("I132").Value = 835716.2
$ws.Range("K132").Value = 2507148.6
$ws.Range("M132").Value = -2504618.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 338.22223
$ws.Range("J80").Value = 408.75
$ws.Range("L80").Value = 408.75
$ws.Range("N80").Value = -2404.75
$ws.Range("H83").Value = 338.22223
$ws.Range("J83").Value = 408.75
$ws.Range("L83").Value = 2043.75
$ws.Range("N83").Value = -12027.75
$ws.Range("H94").Value = 2306.5715
$ws.Range("I94").Value = 2935.375
$ws.Range("J94").Value = 1468.1666
$ws.Range("K94").Value = 2935.375
$ws.Range("L94").Value = 1468.1666
$ws.Range("M94").Value = -2484.375
$ws.Range("N94").Value = -2370.1666
$ws.Range("H107").Value = 1275
$ws.Range("I107").Value = 981.63635
$ws.Range("K107").Value = 981.63635
$ws.Range("M107").Value = 938.36365
$ws.Range("H134").Value = 11921333
$ws.Range("I134").Value = 11198.7
$ws.Range("J134").Value = 41696668
$ws.Range("K134").Value = 33596.10000000001
$ws.Range("L134").Value = 125090004
$ws.Range("M134").Value = -31061.10000000001
$ws.Range("N134").Value = -125095074

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4277446.5
$ws.Range("I31").Value = 6948512
$ws.Range("J31").Value = 3741.8
$ws.Range("K31").Value = 6948512
$ws.Range("L31").Value = 3741.8
$ws.Range("M31").Value = -6948217
$ws.Range("N31").Value = -4331.8
$ws.Range("H34").Value = 4277446.5
$ws.Range("I34").Value = 6948512
$ws.Range("J34").Value = 3741.8
$ws.Range("K34").Value = 6948512
$ws.Range("L34").Value = 3741.8
$ws.Range("M34").Value = -6948310
$ws.Range("N34").Value = -4145.8
$ws.Range("H86").Value = 10413.9375
$ws.Range("J86").Value = 5823
$ws.Range("L86").Value = 5823
$ws.Range("N86").Value = -8069
$ws.Range("H89").Value = 10413.9375
$ws.Range("J89").Value = 5823
$ws.Range("L89").Value = 29115
$ws.Range("N89").Value = -40347
$ws.Range("H107").Value = 375.7857
$ws.Range("I107").Value = 387.6
$ws.Range("J107").Value = 346.25
$ws.Range("K107").Value = 387.6
$ws.Range("L107").Value = 346.25
$ws.Range("M107").Value = 1532.4
$ws.Range("N107").Value = -4186.25
$ws.Range("H122").Value = 19320
$ws.Range("I122").Value = 4876.75
$ws.Range("K122").Value = 14630.25
$ws.Range("M122").Value = -12180.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 200000050
$ws.Range("I11").Value = 58
$ws.Range("K11").Value = 174
$ws.Range("M11").Value = -34
$ws.Range("H116").Value = 3964.5
$ws.Range("J116").Value = 4000
$ws.Range("L116").Value = 12000
$ws.Range("N116").Value = -18884
$ws.Range("H130").Value = 13375.667
$ws.Range("I130").Value = 9999
$ws.Range("J130").Value = 14051
$ws.Range("K130").Value = 29997
$ws.Range("L130").Value = 42153
$ws.Range("M130").Value = -24977
$ws.Range("N130").Value = -52193
$ws.Range("H140").Value = 2412.1667
$ws.Range("I140").Value = 1827.2222
$ws.Range("K140").Value = 5481.6666
$ws.Range("M140").Value = -301.6665999999996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 14999.667
$ws.Range("H50").Value = 14999.667
$ws.Range("H57").Value = 15908.546
$ws.Range("J57").Value = 29999
$ws.Range("L57").Value = 29999
$ws.Range("N57").Value = -31639
$ws.Range("H80").Value = 6917.6665
$ws.Range("I80").Value = 4901.7144
$ws.Range("J80").Value = 9740
$ws.Range("K80").Value = 4901.7144
$ws.Range("L80").Value = 9740
$ws.Range("M80").Value = -3903.7144
$ws.Range("N80").Value = -11736
$ws.Range("H83").Value = 6917.6665
$ws.Range("I83").Value = 4901.7144
$ws.Range("J83").Value = 9740
$ws.Range("K83").Value = 24508.572
$ws.Range("L83").Value = 48700
$ws.Range("M83").Value = -19516.572
$ws.Range("N83").Value = -58684
$ws.Range("H122").Value = 66582.72
$ws.Range("I122").Value = 88576.234
$ws.Range("K122").Value = 265728.702
$ws.Range("M122").Value = -263278.702

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2003.3334
$ws.Range("I22").Value = 599
$ws.Range("J22").Value = 2514
$ws.Range("K22").Value = 599
$ws.Range("L22").Value = 2514
$ws.Range("M22").Value = -304
$ws.Range("N22").Value = -3104
$ws.Range("H27").Value = 2003.3334
$ws.Range("I27").Value = 599
$ws.Range("J27").Value = 2514
$ws.Range("K27").Value = 599
$ws.Range("L27").Value = 2514
$ws.Range("M27").Value = -492
$ws.Range("N27").Value = -2728
$ws.Range("H40").Value = 4193.5
$ws.Range("I40").Value = 3832.6667
$ws.Range("J40").Value = 5997.6665
$ws.Range("K40").Value = 3832.6667
$ws.Range("L40").Value = 5997.6665
$ws.Range("M40").Value = -3696.6667
$ws.Range("N40").Value = -6269.6665
$ws.Range("H46").Value = 3911.5
$ws.Range("I46").Value = 521.3333
$ws.Range("J46").Value = 4693.846
$ws.Range("K46").Value = 521.3333
$ws.Range("L46").Value = 4693.846
$ws.Range("M46").Value = -333.3333
$ws.Range("N46").Value = -5069.846
$ws.Range("H61").Value = 8789.291999999999
$ws.Range("I61").Value = 8097.4287
$ws.Range("K61").Value = 8097.4287
$ws.Range("M61").Value = -7895.4287
$ws.Range("H68").Value = 1750.1538
$ws.Range("I68").Value = 1726.3043
$ws.Range("K68").Value = 1726.3043
$ws.Range("M68").Value = -977.3043
$ws.Range("H71").Value = 1750.1538
$ws.Range("I71").Value = 1726.3043
$ws.Range("K71").Value = 8631.521500000001
$ws.Range("M71").Value = -4887.521500000001
$ws.Range("H113").Value = 8789.291999999999
$ws.Range("I113").Value = 8097.4287
$ws.Range("K113").Value = 8097.4287
$ws.Range("M113").Value = -5927.4287
$ws.Range("H132").Value = 10625249
$ws.Range("I132").Value = 12985560
$ws.Range("K132").Value = 38956680
$ws.Range("M132").Value = -38954150

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9466.9
$ws.Range("I62").Value = 5860
$ws.Range("K62").Value = 5860
$ws.Range("M62").Value = -5236
$ws.Range("H65").Value = 9466.9
$ws.Range("I65").Value = 5860
$ws.Range("K65").Value = 29300
$ws.Range("M65").Value = -26180
$ws.Range("H113").Value = 664.5111000000001
$ws.Range("I113").Value = 658.9666999999999
$ws.Range("K113").Value = 1976.9001
$ws.Range("M113").Value = 193.0999000000002
$ws.Range("H122").Value = 85771.64
$ws.Range("I122").Value = 7387.8887
$ws.Range("J122").Value = 226862.4
$ws.Range("K122").Value = 22163.6661
$ws.Range("L122").Value = 680587.2
$ws.Range("M122").Value = -19713.6661
$ws.Range("N122").Value = -685487.2
